$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: make bold (adds new font + cellXfs entry) ---
$ws.Range("A1:G1").Font.Bold = $true

# --- Row 2 (Qantas) updated figures ---
$ws.Range("B2").Value = 140
$ws.Range("C2").Value = 11.5
$ws.Range("D2").Value = 0.9
$ws.Range("E2").Value = "Bangkok"
$ws.Range("F2").Value = "Melbourne"
$ws.Range("G2").Value = 45014

# --- Row 3 (now Malaysian Air) updated figures ---
$ws.Range("A3").Value = "Malaysian Air"
$ws.Range("B3").Value = 95
$ws.Range("C3").Value = 10.8
$ws.Range("D3").Value = 1.05
$ws.Range("E3").Value = "Kuala Lumpur"
$ws.Range("F3").Value = "Perth"
$ws.Range("G3").Value = 45000

# --- Row 4 (now Air China) updated figures ---
$ws.Range("A4").Value = "Air China"
$ws.Range("B4").Value = 155
$ws.Range("C4").Value = 11.2
$ws.Range("D4").Value = 1.25
$ws.Range("E4").Value = "Shanghai"
$ws.Range("F4").Value = "Melbourne"
$ws.Range("G4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("G4").Value = 45017

# --- Row 5 (new data: Delta Air) ---
$ws.Range("A5").Value = "Delta Air"
$ws.Range("B5").Value = 175
$ws.Range("C5").Value = 10.5
$ws.Range("D5").Value = 1.05
$ws.Range("E5").Value = "Los Angeles"
$ws.Range("G5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("G5").Value = 45046

# --- Row 6 (new data: Singapore Airline, different route) ---
$ws.Range("A6").Value = "Singapore Airline"
$ws.Range("B6").Value = 145
$ws.Range("C6").Value = 9.5
$ws.Range("D6").Value = 1.25
$ws.Range("E6").Value = "Singapore"
$ws.Range("F6").Value = "Perth"
$ws.Range("G6").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("G6").Value = 45031

# Filled in after row 6 (matches original authoring order)
$ws.Range("F5").Value = "Adelaide"

# --- Move the active selection to match the saved view ---
$ws.Range("F9").Select()
